$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (columns H, I, J, K, M, O, P)
# Excel's ColumnWidth property is offset from the stored XML width by ~0.8333333333333334,
# so subtract that offset to land on the target stored widths of 21, 21, 22, 21, 19, 20, 20.
$ws.Columns("H:H").ColumnWidth = 20.166666666666668
$ws.Columns("I:I").ColumnWidth = 20.166666666666668
$ws.Columns("J:J").ColumnWidth = 21.166666666666668
$ws.Columns("K:K").ColumnWidth = 20.166666666666668
$ws.Columns("M:M").ColumnWidth = 18.166666666666668
$ws.Columns("O:O").ColumnWidth = 19.166666666666668
$ws.Columns("P:P").ColumnWidth = 19.166666666666668

# Update cell values across rows 17-105
$ws.Range("D17").Value = 3.6
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 4
$ws.Range("O17").Value = -2.392000000000001
$ws.Range("E18").Value = 10.3
$ws.Range("J18").Value = 10.5
$ws.Range("K18").Value = 10.1
$ws.Range("P18").Value = 48.8
$ws.Range("D19").Value = -8.75
$ws.Range("H19").Value = -8.949999999999999
$ws.Range("I19").Value = -8.550000000000001
$ws.Range("O19").Value = -14.742
$ws.Range("D20").Value = 6.4
$ws.Range("H20").Value = 4.2
$ws.Range("I20").Value = 8.600000000000001
$ws.Range("O20").Value = 44.9
$ws.Range("D21").Value = -9.35
$ws.Range("H21").Value = -9.549999999999999
$ws.Range("I21").Value = -9.15
$ws.Range("O21").Value = 78.65000000000001
$ws.Range("D22").Value = -9.35
$ws.Range("H22").Value = -9.549999999999999
$ws.Range("I22").Value = -9.15
$ws.Range("O22").Value = -97.34999999999999
$ws.Range("D24").Value = 9.25
$ws.Range("H24").Value = 8.800000000000001
$ws.Range("I24").Value = 9.699999999999999
$ws.Range("O24").Value = 109.15
$ws.Range("D26").Value = 7.05
$ws.Range("H26").Value = 5.5
$ws.Range("I26").Value = 8.6
$ws.Range("O26").Value = 46.113
$ws.Range("E27").Value = 12.55
$ws.Range("G27").Value = 0.5
$ws.Range("K27").Value = 12.3
$ws.Range("N27").Value = -50.624
$ws.Range("P27").Value = 63.17400000000001
$ws.Range("W27").Value = -0.062
$ws.Range("D28").Value = -10.6
$ws.Range("F28").Value = 0.5
$ws.Range("H28").Value = -10.85
$ws.Range("I28").Value = -10.35
$ws.Range("M28").Value = -50.624
$ws.Range("O28").Value = 40.024
$ws.Range("V28").Value = -0.062
$ws.Range("D97").Value = 10.35
$ws.Range("H97").Value = 9.9
$ws.Range("I97").Value = 10.8
$ws.Range("O97").Value = 110.25
$ws.Range("D98").Value = 3.6
$ws.Range("H98").Value = 3.2
$ws.Range("I98").Value = 4
$ws.Range("O98").Value = -2.392000000000001
$ws.Range("D99").Value = 2.8
$ws.Range("H99").Value = 2.6
$ws.Range("I99").Value = 3
$ws.Range("O99").Value = 41.59999999999999
$ws.Range("D100").Value = 2.8
$ws.Range("H100").Value = 2.6
$ws.Range("I100").Value = 3
$ws.Range("O100").Value = 41.59999999999999
$ws.Range("D103").Value = -8.15
$ws.Range("E103").Value = 0.8
$ws.Range("H103").Value = -8.35
$ws.Range("I103").Value = -7.95
$ws.Range("J103").Value = 12
$ws.Range("K103").Value = -10.4
$ws.Range("O103").Value = -14.142
$ws.Range("P103").Value = 28.296
$ws.Range("D104").Value = -9.949999999999999
$ws.Range("H104").Value = -10.15
$ws.Range("I104").Value = -9.75
$ws.Range("O104").Value = 78.05
$ws.Range("D105").Value = -9.949999999999999
$ws.Range("H105").Value = -10.15
$ws.Range("I105").Value = -9.75
$ws.Range("O105").Value = -97.95

